$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(3, 65854, 10124, 11002),
    @(4, 38942, 3552, 4826),
    @(5, 102616, 9826, 9281),
    @(6, 1897, 1002, 285),
    @(7, 68832, 11721, 8514),
    @(8, 8012, 1596, 1716),
    @(9, 8168, 1751, 1136),
    @(10, 3920, 599, 320),
    @(11, 1505, 429, 4),
    @(12, 0, 0, 0),
    @(13, 1332, 289, 415),
    @(14, 3817, 1466, 1400),
    @(15, 7013, 2449, 1402),
    @(16, 4798, 2596, 589),
    @(17, 3870, 1259, 477),
    @(18, 25678, 3599, 4791),
    @(19, 2018, 944, 523),
    @(20, 25762, 3540, 3993),
    @(21, 527, 603, 165),
    @(22, 25597, 2869, 4394),
    @(23, 1791, 974, 289),
    @(24, 28578, 2859, 5294),
    @(25, 114043, 9676, 12877),
    @(26, 8688, 3438, 1226),
    @(27, 0, 0, 0),
    @(28, 7655, 1778, 1884),
    @(29, 3493, 813, 733),
    @(30, 21743, 3794, 4339),
    @(31, 641, 149, 386),
    @(32, 4209, 2421, 439),
    @(33, 20507, 4853, 4165),
    @(34, 16191, 4691, 2972),
    @(35, 7260, 934, 1489),
    @(36, 82264, 9370, 8060),
    @(37, 11617, 4312, 1577),
    @(38, 39553, 2815, 3773),
    @(39, 1753, 1548, 277),
    @(40, 2026, 756, 803),
    @(41, 3772, 829, 156),
    @(42, 16141, 807, 436),
    @(43, 388, 290, 73),
    @(44, 1218, 236, 89),
    @(45, 0, 0, 0),
    @(46, 4838, 1569, 595),
    @(47, 20451, 4986, 3897),
    @(48, 48486, 4920, 7246),
    @(49, 22399, 4957, 2040),
    @(50, 18167, 2241, 3969),
    @(51, 47776, 4567, 7538),
    @(52, 7629, 1432, 1865),
    @(53, 16081, 3197, 2700),
    @(54, 3272, 1955, 1364),
    @(55, 3128, 1996, 219),
    @(56, 7655, 1597, 3133),
    @(57, 16960, 6800, 3806),
    @(58, 21519, 2162, 754),
    @(59, 993571, 150974, 141136)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
